$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.225.97"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").Value = "3.313.37"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "'526.21"
$ws.Range("E5").Value = "  -2.29%  "

$ws.Range("D6").Value = "'172.43"
$ws.Range("E6").Value = "  -6.94%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.585"
$ws.Range("E7").Value = "  -3.64%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.315.70"
$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").Value = "'0.603"
$ws.Range("E10").Value = "  -3.74%  "

$ws.Range("D11").Value = "'52.83"
$ws.Range("E11").Value = "  -14.03%  "

$ws.Range("E12").Value = "  -1.77%  "

$ws.Range("D13").Value = "'0.0000257"
$ws.Range("E13").Value = "  -3.29%  "

$ws.Range("D14").Value = "'8.92"
$ws.Range("E14").Value = "  -3.32%  "

$ws.Range("D15").Value = "3.832.87"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.304.99"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.117"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").Value = "64.078.05"
$ws.Range("E18").Value = "  -2.04%  "

$ws.Range("D19").Value = "'17.41"
$ws.Range("E19").Value = "  -2.81%  "

$ws.Range("D20").Value = "'11.17"
$ws.Range("E20").Value = "  -1.02%  "

$ws.Range("D21").Value = "'0.954"

$ws.Range("D22").Value = "'378.61"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").Value = "'4.17"
$ws.Range("E23").Value = "  +6.94%  "

$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").Value = "'11.10"
$ws.Range("E25").Value = "  -3.10%  "

$ws.Range("D26").Value = "'3.69"
$ws.Range("E26").Value = "  -5.08%  "

$ws.Range("D27").Value = "'6.15"
$ws.Range("E27").Value = "  +1.44%  "

$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").Value = "'11.22"
$ws.Range("E29").Value = "  -4.39%  "

$ws.Range("D30").Value = "'8.12"
$ws.Range("E30").Value = "  -5.33%  "

$ws.Range("D31").Value = "'28.71"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D32").Value = "'626.41"
$ws.Range("E32").Value = "  -4.82%  "

$ws.Range("D33").Value = "'6.56"
$ws.Range("E33").Value = "  -4.77%  "

$ws.Range("D34").Value = "'11.17"
$ws.Range("E34").Value = "  -2.36%  "

$ws.Range("E35").Value = "  -1.89%  "

$ws.Range("D36").Value = "'56.89"
$ws.Range("E36").Value = "  -4.97%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").Value = "'35.94"
$ws.Range("E38").Value = "  -3.27%  "

$ws.Range("D39").Value = "'0.377"
$ws.Range("E39").Value = "  -5.42%  "

$ws.Range("D40").Value = "0.0₃0744"
$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("D42").Value = "'3.18"
$ws.Range("E42").Value = "  +8.37%  "

$ws.Range("D43").Value = "'0.124"
$ws.Range("E43").Value = "  -3.65%  "

$ws.Range("D44").Value = "'2.60"
$ws.Range("E44").Value = "  +2.29%  "

$ws.Range("D45").Value = "2.876.85"
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").Value = "'0.0396"
$ws.Range("E47").Value = "  -2.63%  "

$ws.Range("E48").Value = "  -5.23%  "

$ws.Range("D49").Value = "'3.06"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("D50").Value = "'137.47"
$ws.Range("E50").Value = "  +1.29%  "

$ws.Range("E51").Value = "  -2.46%  "

